$d = $word.ActiveDocument

$replacements = @(
    @("2023-08-27 Sunday", "2023-08-28 Monday"),
    @("47÷3=", "41÷8="),
    @("53÷3=", "97÷7="),
    @("29÷5=", "37÷2="),
    @("87÷6=", "49÷5="),
    @("23÷9=", "15÷4="),
    @("95÷7=", "55÷4="),
    @("75÷2=", "50÷7="),
    @("88÷8=", "48÷2="),
    @("79÷6=", "53÷4="),
    @("32÷7=", "55÷8="),
    @("83÷9=", "33÷3="),
    @("33÷5=", "93÷3="),
    @("91÷5=", "60÷9="),
    @("23÷4=", "42÷2="),
    @("38÷8=", "11÷6="),
    @("90÷6=", "80÷5="),
    @("60÷7=", "16÷2="),
    @("91÷3=", "63÷8="),
    @("37÷6=", "90÷3="),
    @("94÷4=", "37÷7="),
    @("39÷4=", "90÷8="),
    @("59÷7=", "26÷8="),
    @("81÷8=", "92÷6="),
    @("64÷9=", "75÷8="),
    @("35÷7=", "78÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
